# Update gh-pages to output generated at 456a3b4
# Applies the updated "想去人数" (F column) counts across the three sheets
# that carry this data: 展览, 演出, and 全部类型 (the aggregate sheet).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 456
$ws1.Range("F5").Value  = 68
$ws1.Range("F6").Value  = 1
$ws1.Range("F7").Value  = 1301
$ws1.Range("F13").Value = 1087
$ws1.Range("F17").Value = 85
$ws1.Range("F18").Value = 237
$ws1.Range("F19").Value = 1644
$ws1.Range("F20").Value = 605
$ws1.Range("F22").Value = 184
$ws1.Range("F23").Value = 1790
$ws1.Range("F27").Value = 1201
$ws1.Range("F28").Value = 58
$ws1.Range("F31").Value = 1593
$ws1.Range("F33").Value = 110
$ws1.Range("F34").Value = 626
$ws1.Range("F36").Value = 1746
$ws1.Range("F38").Value = 1771
$ws1.Range("F39").Value = 196
$ws1.Range("F41").Value = 833
$ws1.Range("F42").Value = 35
$ws1.Range("F43").Value = 831
$ws1.Range("F44").Value = 783
$ws1.Range("F45").Value = 992
$ws1.Range("F47").Value = 3314

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 787

# --- Sheet "全部类型" (all types, aggregate) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 456
$ws4.Range("F5").Value  = 68
$ws4.Range("F8").Value  = 1301
$ws4.Range("F14").Value = 1087
$ws4.Range("F18").Value = 85
$ws4.Range("F20").Value = 237
$ws4.Range("F21").Value = 1644
$ws4.Range("F22").Value = 605
$ws4.Range("F24").Value = 184
$ws4.Range("F25").Value = 1790
$ws4.Range("F28").Value = 1201
$ws4.Range("F30").Value = 1593
$ws4.Range("F32").Value = 110
$ws4.Range("F33").Value = 787
$ws4.Range("F35").Value = 626
$ws4.Range("F37").Value = 1746
$ws4.Range("F40").Value = 1771
$ws4.Range("F41").Value = 833
$ws4.Range("F42").Value = 831
$ws4.Range("F43").Value = 783
$ws4.Range("F44").Value = 992
$ws4.Range("F48").Value = 3314
